# Update TPM-derived statistics for the Il6-Il6st LR-pairs sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 25,16
$data[0,0] = 3.0
$data[0,1] = 1.0
$data[0,2] = 8.651718
$data[0,3] = 25.955154
$data[0,4] = 0.09424769069200405
$data[0,5] = 0.09527060788475794
$data[0,6] = 3.0
$data[0,7] = 1.0
$data[0,8] = 55.783591
$data[0,9] = 167.350773
$data[0,10] = 0.2332214199005771
$data[0,11] = 0.2394371967339281
$data[0,12] = 482.6238983593381
$data[0,13] = 4343.615085234042
$data[0,14] = 0.02198058024553959
$data[0,15] = 0.02281132728306371
$data[1,0] = 3.0
$data[1,1] = 1.0
$data[1,2] = 8.651718
$data[1,3] = 25.955154
$data[1,4] = 0.09424769069200405
$data[1,5] = 0.09527060788475794
$data[1,6] = 3.0
$data[1,7] = 1.0
$data[1,8] = 126.7095336666667
$data[1,9] = 380.128601
$data[1,10] = 0.5297503589663128
$data[1,11] = 0.5438691736537713
$data[1,12] = 1096.255153195506
$data[1,13] = 9866.296378759554
$data[1,14] = 0.04992774797583516
$data[1,15] = 0.05181474678377577
$data[2,0] = 3.0
$data[2,1] = 1.0
$data[2,2] = 8.651718
$data[2,3] = 25.955154
$data[2,4] = 0.09424769069200405
$data[2,5] = 0.09527060788475794
$data[2,6] = 3.0
$data[2,7] = 1.0
$data[2,8] = 23.03749833333333
$data[2,9] = 69.112495
$data[2,10] = 0.09631574403765399
$data[2,11] = 0.09888273454277752
$data[2,12] = 199.31393900547
$data[2,13] = 1793.82545104923
$data[2,14] = 0.009077536452831046
$data[2,15] = 0.009420618229197566
$data[3,0] = 3.0
$data[3,1] = 1.0
$data[3,2] = 8.651718
$data[3,3] = 25.955154
$data[3,4] = 0.09424769069200405
$data[3,5] = 0.09527060788475794
$data[3,6] = 2.0
$data[3,7] = 1.0
$data[3,8] = 18.627865
$data[3,9] = 37.25573
$data[3,10] = 0.07787983970082285
$data[3,11] = 0.05330365312071852
$data[3,12] = 161.16303492207
$data[3,13] = 966.97820953242
$data[3,14] = 0.007339995043266009
$data[3,15] = 0.005078271435289128
$data[4,0] = 3.0
$data[4,1] = 1.0
$data[4,2] = 8.651718
$data[4,3] = 25.955154
$data[4,4] = 0.09424769069200405
$data[4,5] = 0.09527060788475794
$data[4,6] = 3.0
$data[4,7] = 1.0
$data[4,8] = 15.028766
$data[4,9] = 45.086298
$data[4,10] = 0.06283263739463307
$data[4,11] = 0.06450724194880479
$data[4,12] = 130.024645319988
$data[4,13] = 1170.221807879892
$data[4,14] = 0.005921830974532225
$data[4,15] = 0.00614564415343179
$data[5,0] = 3.0
$data[5,1] = 1.0
$data[5,2] = 76.28028133333333
$data[5,3] = 228.840844
$data[5,4] = 0.8309610138706613
$data[5,5] = 0.8399798481928121
$data[5,6] = 3.0
$data[5,7] = 1.0
$data[5,8] = 55.783591
$data[5,9] = 167.350773
$data[5,10] = 0.2332214199005771
$data[5,11] = 0.2394371967339281
$data[5,12] = 4255.188015263601
$data[5,13] = 38296.69213737241
$data[5,14] = 0.1937979075369388
$data[5,15] = 0.2011224201642774
$data[6,0] = 3.0
$data[6,1] = 1.0
$data[6,2] = 76.28028133333333
$data[6,3] = 228.840844
$data[6,4] = 0.8309610138706613
$data[6,5] = 0.8399798481928121
$data[6,6] = 3.0
$data[6,7] = 1.0
$data[6,8] = 126.7095336666667
$data[6,9] = 380.128601
$data[6,10] = 0.5297503589663128
$data[6,11] = 0.5438691736537713
$data[6,12] = 9665.438875708805
$data[6,13] = 86988.94988137925
$data[6,14] = 0.4402018953849941
$data[6,15] = 0.456839145922445
$data[7,0] = 3.0
$data[7,1] = 1.0
$data[7,2] = 76.28028133333333
$data[7,3] = 228.840844
$data[7,4] = 0.8309610138706613
$data[7,5] = 0.8399798481928121
$data[7,6] = 3.0
$data[7,7] = 1.0
$data[7,8] = 23.03749833333333
$data[7,9] = 69.112495
$data[7,10] = 0.09631574403765399
$data[7,11] = 0.09888273454277752
$data[7,12] = 1757.306854082864
$data[7,13] = 15815.76168674578
$data[7,14] = 0.08003462831723607
$data[7,15] = 0.0830595043501324
$data[8,0] = 3.0
$data[8,1] = 1.0
$data[8,2] = 76.28028133333333
$data[8,3] = 228.840844
$data[8,4] = 0.8309610138706613
$data[8,5] = 0.8399798481928121
$data[8,6] = 2.0
$data[8,7] = 1.0
$data[8,8] = 18.627865
$data[8,9] = 37.25573
$data[8,10] = 0.07787983970082285
$data[8,11] = 0.05330365312071852
$data[8,12] = 1420.938782839353
$data[8,13] = 8525.63269703612
$data[8,14] = 0.06471511055788035
$data[8,15] = 0.04477399445646345
$data[9,0] = 3.0
$data[9,1] = 1.0
$data[9,2] = 76.28028133333333
$data[9,3] = 228.840844
$data[9,4] = 0.8309610138706613
$data[9,5] = 0.8399798481928121
$data[9,6] = 3.0
$data[9,7] = 1.0
$data[9,8] = 15.028766
$data[9,9] = 45.086298
$data[9,10] = 0.06283263739463307
$data[9,11] = 0.06450724194880479
$data[9,12] = 1146.398498572835
$data[9,13] = 10317.58648715551
$data[9,14] = 0.05221147207361193
$data[9,15] = 0.05418478329949405
$data[10,0] = 3.0
$data[10,1] = 1.0
$data[10,2] = 3.806801666666667
$data[10,3] = 11.420405
$data[10,4] = 0.04146948224685611
$data[10,5] = 0.04191957122042617
$data[10,6] = 3.0
$data[10,7] = 1.0
$data[10,8] = 55.783591
$data[10,9] = 167.350773
$data[10,10] = 0.2332214199005771
$data[10,11] = 0.2394371967339281
$data[10,12] = 212.3570671914517
$data[10,13] = 1911.213604723065
$data[10,14] = 0.009671571532153558
$data[10,15] = 0.01003710462130709
$data[11,0] = 3.0
$data[11,1] = 1.0
$data[11,2] = 3.806801666666667
$data[11,3] = 11.420405
$data[11,4] = 0.04146948224685611
$data[11,5] = 0.04191957122042617
$data[11,6] = 3.0
$data[11,7] = 1.0
$data[11,8] = 126.7095336666667
$data[11,9] = 380.128601
$data[11,10] = 0.5297503589663128
$data[11,11] = 0.5438691736537713
$data[11,12] = 482.3580639448228
$data[11,13] = 4341.222575503405
$data[11,14] = 0.02196847310641916
$data[11,15] = 0.02279876255957359
$data[12,0] = 3.0
$data[12,1] = 1.0
$data[12,2] = 3.806801666666667
$data[12,3] = 11.420405
$data[12,4] = 0.04146948224685611
$data[12,5] = 0.04191957122042617
$data[12,6] = 3.0
$data[12,7] = 1.0
$data[12,8] = 23.03749833333333
$data[12,9] = 69.112495
$data[12,10] = 0.09631574403765399
$data[12,11] = 0.09888273454277752
$data[12,12] = 87.69918705116389
$data[12,13] = 789.292683460475
$data[12,14] = 0.003994164037462229
$data[12,15] = 0.004145121833136456
$data[13,0] = 3.0
$data[13,1] = 1.0
$data[13,2] = 3.806801666666667
$data[13,3] = 11.420405
$data[13,4] = 0.04146948224685611
$data[13,5] = 0.04191957122042617
$data[13,6] = 2.0
$data[13,7] = 1.0
$data[13,8] = 18.627865
$data[13,9] = 37.25573
$data[13,10] = 0.07787983970082285
$data[13,11] = 0.05330365312071852
$data[13,12] = 70.91258752844168
$data[13,13] = 425.47552517065
$data[13,14] = 0.003229636629861273
$data[13,15] = 0.002234466283302851
$data[14,0] = 3.0
$data[14,1] = 1.0
$data[14,2] = 3.806801666666667
$data[14,3] = 11.420405
$data[14,4] = 0.04146948224685611
$data[14,5] = 0.04191957122042617
$data[14,6] = 3.0
$data[14,7] = 1.0
$data[14,8] = 15.028766
$data[14,9] = 45.086298
$data[14,10] = 0.06283263739463307
$data[14,11] = 0.06450724194880479
$data[14,12] = 57.21153145674333
$data[14,13] = 514.90378311069
$data[14,14] = 0.002605636940959884
$data[14,15] = 0.002704115923106185
$data[15,0] = 2.0
$data[15,1] = 1.0
$data[15,2] = 2.956885
$data[15,3] = 5.91377
$data[15,4] = 0.03221089532643417
$data[15,5] = 0.02170699749231482
$data[15,6] = 3.0
$data[15,7] = 1.0
$data[15,8] = 55.783591
$data[15,9] = 167.350773
$data[15,10] = 0.2332214199005771
$data[15,11] = 0.2394371967339281
$data[15,12] = 164.945663474035
$data[15,13] = 989.67398084421
$data[15,14] = 0.00751227074429984
$data[15,15] = 0.005197462629070266
$data[16,0] = 2.0
$data[16,1] = 1.0
$data[16,2] = 2.956885
$data[16,3] = 5.91377
$data[16,4] = 0.03221089532643417
$data[16,5] = 0.02170699749231482
$data[16,6] = 3.0
$data[16,7] = 1.0
$data[16,8] = 126.7095336666667
$data[16,9] = 380.128601
$data[16,10] = 0.5297503589663128
$data[16,11] = 0.5438691736537713
$data[16,12] = 374.6655194559617
$data[16,13] = 2247.99311673577
$data[16,14] = 0.01706373336180483
$data[16,15] = 0.01180576678864974
$data[17,0] = 2.0
$data[17,1] = 1.0
$data[17,2] = 2.956885
$data[17,3] = 5.91377
$data[17,4] = 0.03221089532643417
$data[17,5] = 0.02170699749231482
$data[17,6] = 3.0
$data[17,7] = 1.0
$data[17,8] = 23.03749833333333
$data[17,9] = 69.112495
$data[17,10] = 0.09631574403765399
$data[17,11] = 0.09888273454277752
$data[17,12] = 68.11923325935832
$data[17,13] = 408.7153995561499
$data[17,14] = 0.003102416349484498
$data[17,15] = 0.002146447270753303
$data[18,0] = 2.0
$data[18,1] = 1.0
$data[18,2] = 2.956885
$data[18,3] = 5.91377
$data[18,4] = 0.03221089532643417
$data[18,5] = 0.02170699749231482
$data[18,6] = 2.0
$data[18,7] = 1.0
$data[18,8] = 18.627865
$data[18,9] = 37.25573
$data[18,10] = 0.07787983970082285
$data[18,11] = 0.05330365312071852
$data[18,12] = 55.08045460052499
$data[18,13] = 220.3218184021
$data[18,14] = 0.002508579364642677
$data[18,15] = 0.001157062264622656
$data[19,0] = 2.0
$data[19,1] = 1.0
$data[19,2] = 2.956885
$data[19,3] = 5.91377
$data[19,4] = 0.03221089532643417
$data[19,5] = 0.02170699749231482
$data[19,6] = 3.0
$data[19,7] = 1.0
$data[19,8] = 15.028766
$data[19,9] = 45.086298
$data[19,10] = 0.06283263739463307
$data[19,11] = 0.06450724194880479
$data[19,12] = 44.43833275391
$data[19,13] = 266.62999652346
$data[19,14] = 0.002023895506202319
$data[19,15] = 0.001400258539218851
$data[20,0] = 1.0
$data[20,1] = 0.3333333333333333
$data[20,2] = 0.1019796666666667
$data[20,3] = 0.305939
$data[20,4] = 0.001110917864044306
$data[20,5] = 0.001122975209688795
$data[20,6] = 3.0
$data[20,7] = 1.0
$data[20,8] = 55.783591
$data[20,9] = 167.350773
$data[20,10] = 0.2332214199005771
$data[20,11] = 0.2394371967339281
$data[20,12] = 5.688792015649668
$data[20,13] = 51.199128140847
$data[20,14] = 0.0002590898416453293
$data[20,15] = 0.0002688820362095801
$data[21,0] = 1.0
$data[21,1] = 0.3333333333333333
$data[21,2] = 0.1019796666666667
$data[21,3] = 0.305939
$data[21,4] = 0.001110917864044306
$data[21,5] = 0.001122975209688795
$data[21,6] = 3.0
$data[21,7] = 1.0
$data[21,8] = 126.7095336666667
$data[21,9] = 380.128601
$data[21,10] = 0.5297503589663128
$data[21,11] = 0.5438691736537713
$data[21,12] = 12.92179600681545
$data[21,13] = 116.296164061339
$data[21,14] = 0.0005885091372595607
$data[21,15] = 0.0006107515993271154
$data[22,0] = 1.0
$data[22,1] = 0.3333333333333333
$data[22,2] = 0.1019796666666667
$data[22,3] = 0.305939
$data[22,4] = 0.001110917864044306
$data[22,5] = 0.001122975209688795
$data[22,6] = 3.0
$data[22,7] = 1.0
$data[22,8] = 23.03749833333333
$data[22,9] = 69.112495
$data[22,10] = 0.09631574403765399
$data[22,11] = 0.09888273454277752
$data[22,12] = 2.349356400867222
$data[22,13] = 21.144207607805
$data[22,14] = 0.0001069988806401487
$data[22,15] = 0.000111042859557777
$data[23,0] = 1.0
$data[23,1] = 0.3333333333333333
$data[23,2] = 0.1019796666666667
$data[23,3] = 0.305939
$data[23,4] = 0.001110917864044306
$data[23,5] = 0.001122975209688795
$data[23,6] = 2.0
$data[23,7] = 1.0
$data[23,8] = 18.627865
$data[23,9] = 37.25573
$data[23,10] = 0.07787983970082285
$data[23,11] = 0.05330365312071852
$data[23,12] = 1.899663463411667
$data[23,13] = 11.39798078047
$data[23,14] = 0.00008651810517255106
$data[23,15] = 0.00005985868104041766
$data[24,0] = 1.0
$data[24,1] = 0.3333333333333333
$data[24,2] = 0.1019796666666667
$data[24,3] = 0.305939
$data[24,4] = 0.001110917864044306
$data[24,5] = 0.001122975209688795
$data[24,6] = 3.0
$data[24,7] = 1.0
$data[24,8] = 15.028766
$data[24,9] = 45.086298
$data[24,10] = 0.06283263739463307
$data[24,11] = 0.06450724194880479
$data[24,12] = 1.532628547091333
$data[24,13] = 13.793656923822
$data[24,14] = 0.00006980189932671617
$data[24,15] = 0.00007244003355390488

$ws.Range("E2:T26").Value2 = $data
